# Add 2022-Q3 data:
#   - a new worksheet "2022-Q3" with the per-fund holdings detail, inserted
#     right after "总计" (and therefore right before "2022-Q2"), and
#   - a new summary row on "总计" (row 2), pushing the existing rows down.

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item("总计")
$q2Sheet = $wb.Worksheets.Item("2022-Q2")

# ---------------------------------------------------------------------
# 1) Build the new "2022-Q3" detail sheet by duplicating the "2022-Q2"
#    sheet (so it inherits the same column widths / header & cell
#    styles), then trim it down to the 4 data rows we need and
#    overwrite the values.
# ---------------------------------------------------------------------
$q2Sheet.Copy($null, $totalSheet)
$q3Sheet = $wb.Worksheets.Item($totalSheet.Index + 1)
$q3Sheet.Name = "2022-Q3"

# The template ("2022-Q2") has 8 data rows (rows 2-9); 2022-Q3 only
# needs 4 (rows 2-5), so drop the extra rows.
$q3Sheet.Rows("6:9").Delete()

# Keep the fund-code-like text columns (B,C,D,E,F,G) as plain text so
# leading zeros / decimal-looking codes are not coerced into numbers.
$q3Sheet.Range("B2:G5").NumberFormat = "@"

$q3Data = @(
    @(0, "673060", "西部利得景瑞灵活配置混合A", "3.97", "93.10", "4.46", "0.1771", 4),
    @(1, "006234", "万家汽车新趋势混合C",       "2.23", "90.68", "2.89", "0.0644", 9),
    @(2, "009258", "西部利得景瑞灵活配置混合C", "1.32", "93.10", "4.46", "0.0589", 4),
    @(3, "006233", "万家汽车新趋势混合A",       "1.93", "90.68", "2.89", "0.0558", 9)
)

for ($i = 0; $i -lt $q3Data.Length; $i++) {
    $row = 2 + $i
    $values = $q3Data[$i]
    $q3Sheet.Cells.Item($row, 1).Value = $values[0]
    $q3Sheet.Cells.Item($row, 2).Value = $values[1]
    $q3Sheet.Cells.Item($row, 3).Value = $values[2]
    $q3Sheet.Cells.Item($row, 4).Value = $values[3]
    $q3Sheet.Cells.Item($row, 5).Value = $values[4]
    $q3Sheet.Cells.Item($row, 6).Value = $values[5]
    $q3Sheet.Cells.Item($row, 7).Value = $values[6]
    $q3Sheet.Cells.Item($row, 8).Value = $values[7]
}

# ---------------------------------------------------------------------
# 2) Insert the new summary row into "总计" (row 2), copy the row
#    formatting from the row below (the old row 2, now shifted to row
#    3) and fill in the 2022-Q3 totals.
# ---------------------------------------------------------------------
$totalSheet.Rows("2:2").Insert()
$totalSheet.Range("A3:D3").Copy($totalSheet.Range("A2:D2"))

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 4
$totalSheet.Range("D2").Value = 0.36

# The "index" column (A) is a simple 0-based row counter; bump every
# row that shifted down by one because of the insert above.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
